# Update survey ratings and the selected cell, per the authored change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the rating values that changed.
$ws.Range("F3").Value = 3
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 4

# The saved selection moved from D5 to F3.
$ws.Range("F3").Select()
